$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2035087719298246
$ws.Range("C2").Value = 0.5368421052631579
$ws.Range("J2").Value = 0.01403508771929825
$ws.Range("P2").Value = 0.1333333333333333
$ws.Range("S2").Value = 0.112280701754386

$ws.Range("B3").Value = 0.00625
$ws.Range("C3").Value = 0.04375
$ws.Range("J3").Value = 0.01875
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.18125

$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("P4").Value = 0.6511627906976745
$ws.Range("S4").Value = 0.3023255813953488

$ws.Range("P5").Value = 0.4285714285714285
$ws.Range("S5").Value = 0.5714285714285714

$ws.Range("B6").Value = 0.04597701149425287
$ws.Range("E6").Value = 0.005747126436781609
$ws.Range("F6").Value = 0.03448275862068965
$ws.Range("J6").Value = 0.3218390804597701
$ws.Range("O6").Value = 0.02298850574712644
$ws.Range("Q6").Value = 0.1494252873563219
$ws.Range("R6").Value = 0.08045977011494253
$ws.Range("S6").Value = 0.3390804597701149

$ws.Range("B7").Value = 0.1125
$ws.Range("D7").Value = 0.01875
$ws.Range("E7").Value = 0.00625
$ws.Range("F7").Value = 0.05
$ws.Range("J7").Value = 0.15625
$ws.Range("O7").Value = 0.01875
$ws.Range("Q7").Value = 0.14375
$ws.Range("R7").Value = 0.06875
$ws.Range("S7").Value = 0.425

$ws.Range("B8").Value = 0.1026392961876833
$ws.Range("D8").Value = 0.01173020527859238
$ws.Range("E8").Value = 0.002932551319648094
$ws.Range("F8").Value = 0.04985337243401759
$ws.Range("J8").Value = 0.09090909090909091
$ws.Range("O8").Value = 0.03225806451612903
$ws.Range("Q8").Value = 0.1642228739002932
$ws.Range("R8").Value = 0.07917888563049853
$ws.Range("S8").Value = 0.4662756598240469

$ws.Range("B9").Value = 0.1096774193548387
$ws.Range("D9").Value = 0.006451612903225806
$ws.Range("E9").Value = 0.01290322580645161
$ws.Range("F9").Value = 0.05806451612903226
$ws.Range("J9").Value = 0.1161290322580645
$ws.Range("O9").Value = 0.01935483870967742
$ws.Range("Q9").Value = 0.1096774193548387
$ws.Range("R9").Value = 0.06451612903225806
$ws.Range("S9").Value = 0.5032258064516129

$ws.Range("B10").Value = 0.1463178294573644
$ws.Range("D10").Value = 0.03488372093023256
$ws.Range("E10").Value = 0.001937984496124031
$ws.Range("F10").Value = 0.07170542635658915
$ws.Range("J10").Value = 0.1007751937984496
$ws.Range("O10").Value = 0.02616279069767442
$ws.Range("Q10").Value = 0.1831395348837209
$ws.Range("R10").Value = 0.06492248062015504
$ws.Range("S10").Value = 0.3701550387596899

$ws.Range("F11").Value = 0.0
$ws.Range("G11").Value = 0.1704545454545454
$ws.Range("J11").Value = 0.0946969696969697
$ws.Range("K11").Value = 0.2045454545454546
$ws.Range("L11").Value = 0.5189393939393939
$ws.Range("S11").Value = 0.01136363636363636

$ws.Range("G12").Value = 0.7214285714285714
$ws.Range("J12").Value = 0.1785714285714286
$ws.Range("K12").Value = 0.007142857142857143
$ws.Range("L12").Value = 0.02857142857142857
$ws.Range("S12").Value = 0.06428571428571428

$ws.Range("G13").Value = 0.5757575757575758
$ws.Range("J13").Value = 0.3636363636363636
$ws.Range("S13").Value = 0.06060606060606061

$ws.Range("F15").Value = 0.01869158878504673
$ws.Range("H15").Value = 0.1308411214953271
$ws.Range("I15").Value = 0.06542056074766354
$ws.Range("J15").Value = 0.3598130841121495
$ws.Range("K15").Value = 0.07009345794392523
$ws.Range("M15").Value = 0.004672897196261682
$ws.Range("O15").Value = 0.05607476635514019
$ws.Range("S15").Value = 0.294392523364486

$ws.Range("F16").Value = 0.01098901098901099
$ws.Range("H16").Value = 0.1208791208791209
$ws.Range("I16").Value = 0.09340659340659341
$ws.Range("J16").Value = 0.4010989010989011
$ws.Range("K16").Value = 0.1098901098901099
$ws.Range("M16").Value = 0.03846153846153846
$ws.Range("O16").Value = 0.07142857142857142
$ws.Range("S16").Value = 0.1538461538461539

$ws.Range("F17").Value = 0.02580645161290323
$ws.Range("H17").Value = 0.1806451612903226
$ws.Range("I17").Value = 0.1096774193548387
$ws.Range("J17").Value = 0.3709677419354839
$ws.Range("K17").Value = 0.0967741935483871
$ws.Range("M17").Value = 0.01612903225806452
$ws.Range("O17").Value = 0.07096774193548387
$ws.Range("S17").Value = 0.1290322580645161

$ws.Range("H18").Value = 0.1317829457364341
$ws.Range("I18").Value = 0.1085271317829457
$ws.Range("J18").Value = 0.4883720930232558
$ws.Range("K18").Value = 0.124031007751938
$ws.Range("O18").Value = 0.05426356589147287
$ws.Range("S18").Value = 0.09302325581395349

$ws.Range("F19").Value = 0.02398523985239853
$ws.Range("H19").Value = 0.202029520295203
$ws.Range("I19").Value = 0.07103321033210332
$ws.Range("J19").Value = 0.3809963099630996
$ws.Range("K19").Value = 0.1171586715867159
$ws.Range("M19").Value = 0.01937269372693727
$ws.Range("N19").Value = 0.0009225092250922509
$ws.Range("O19").Value = 0.07933579335793357
$ws.Range("S19").Value = 0.1051660516605166
